# Updated page creator class to include selects
#
# - Rename the two generated "jultwo" test client rows to the new
#   JulTest naming (keeping B3's first name string as-is per the diff).
# - Move the sheet selection from B3 to B6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "cLNameJulTestA"
$ws.Range("B2").Value = "cFNameJulTestA"
$ws.Range("A3").Value = "cLNameJulTestB"
$ws.Range("B3").Value = "cFNameJultwoB"

$ws.Range("B6").Select()
